$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label text casing updates (shared strings) ---
$ws.Range("C1").Value = "Integral"
$ws.Range("D1").Value = "Time"
$ws.Range("A13").Value = "Avg"

# --- New brute-force run data (B2:D11) ---
$ws.Range("B2").Value = 0.00095089800000000004
$ws.Range("C2").Value = 0.1910075705
$ws.Range("D2").Value = 98.045579169999996

$ws.Range("B3").Value = 0.0010855905569999999
$ws.Range("C3").Value = 0.19309562690000001
$ws.Range("D3").Value = 98.342571530000001

$ws.Range("B4").Value = 0.00098790553029999995
$ws.Range("C4").Value = 0.19272740629999999
$ws.Range("D4").Value = 98.2377939

$ws.Range("B5").Value = 0.0010851493699999999
$ws.Range("C5").Value = 0.194158211
$ws.Range("D5").Value = 99.778155569999996

$ws.Range("B6").Value = 0.0010595598529999999
$ws.Range("C6").Value = 0.19305126789999999
$ws.Range("D6").Value = 103.5504213

$ws.Range("B7").Value = 0.0010962150900000001
$ws.Range("C7").Value = 0.19400646260000001
$ws.Range("D7").Value = 98.210728470000006

$ws.Range("B8").Value = 0.0012973485600000001
$ws.Range("C8").Value = 0.1931227515
$ws.Range("D8").Value = 98.252529440000004

$ws.Range("B9").Value = 0.0010909581759999999
$ws.Range("C9").Value = 0.1936196717
$ws.Range("D9").Value = 99.112405330000001

$ws.Range("B10").Value = 0.001041741636
$ws.Range("C10").Value = 0.19184460710000001
$ws.Range("D10").Value = 98.081835720000001

$ws.Range("B11").Value = 0.0010010040880000001
$ws.Range("C11").Value = 0.1924029300
$ws.Range("D11").Value = 97.994210580000001

# --- Formulas in the AVG / STD summary rows stay the same; just let them recalc ---
$ws.Range("B13").Formula = "=AVERAGE(B2:B11)"
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"
$ws.Range("D13").Formula = "=AVERAGE(D2:D11)"

$ws.Range("B14").Formula = "=_xlfn.STDEV.S(B2:B11)"
$ws.Range("C14").Formula = "=_xlfn.STDEV.S(C2:C11)"
$ws.Range("D14").Formula = "=_xlfn.STDEV.S(D2:D11)"

# --- Whole-sheet default font bumped from 11pt to 12pt ---
$ws.Cells.Font.Size = 12

# --- Row height / column width follow the new font size ---
$ws.Range("A1:D14").RowHeight = 16
$ws.Columns("B").EntireColumn.AutoFit()

# --- Page margins reset to Excel's defaults ---
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# --- Final selection lands on D14 ---
$ws.Range("D14").Select()
